$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("總表")
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
